$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header row 15: drop the "(S)"/"(LO)" qualifiers on the splitter labels
$ws.Range("D15").Value = "Splitter S"
$ws.Range("E15").Value = "Splitter "

# Column A labels (rows 16-21): drop the unit suffixes
$ws.Range("A16").Value = "S Pav"
$ws.Range("A17").Value = "LO Pav"
$ws.Range("A18").Value = "Penalty"
$ws.Range("A20").Value = "S aggr Penalty"
$ws.Range("A21").Value = "LO aggr Penalty"

# Widen column D a bit (was 6.5703125 characters ~= 46px, now ~62px)
$ws.Range("D1").ColumnWidth = 8

# The embedded Visio object is anchored with a fixed right edge in EMUs; widening
# column D pushes that edge further right in cell-relative terms, so pull the
# object back in by the same absolute amount to keep its on-screen size/position.
$shp = $ws.Shapes.Item(1)
$shp.Width = $shp.Width - 12.0008

# Restore the cursor/selection position that was saved with the sheet
$ws.Activate()
$ws.Range("A22").Select()
